# Scheduled-runner refresh of market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# per-crafting-class leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H17").Value = 803
$ws.Range("J17").Value = 759.8889
$ws.Range("L17").Value = 2279.6667
$ws.Range("N17").Value = -2615.6667

# Row 21
$ws.Range("H20").Value = 837
$ws.Range("I20").Value = 837
$ws.Range("K20").Value = 837
$ws.Range("M20").Value = -607

# Row 36
$ws.Range("H35").Value = 837
$ws.Range("I35").Value = 837
$ws.Range("K35").Value = 837
$ws.Range("M35").Value = -458

# Row 114
$ws.Range("H113").Value = 13082.167
$ws.Range("J113").Value = 5874
$ws.Range("L113").Value = 5874
$ws.Range("N113").Value = -12382

# Row 126
$ws.Range("H125").Value = 3246.5715
$ws.Range("J125").Value = 3816.375
$ws.Range("L125").Value = 34347.375
$ws.Range("N125").Value = -39267.375

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H45").Value = 4085.7144
$ws.Range("J45").Value = 4307.6924
$ws.Range("L45").Value = 4307.6924
$ws.Range("N45").Value = -5061.6924

# Row 111
$ws.Range("H110").Value = 516.3333
$ws.Range("I110").Value = 479.6
$ws.Range("K110").Value = 479.6
$ws.Range("M110").Value = 1565.4

# Row 123
$ws.Range("H122").Value = 1687.6923
$ws.Range("I122").Value = 1687.6923
$ws.Range("K122").Value = 5063.0769
$ws.Range("M122").Value = -2613.0769

$ws = $wb.Worksheets.Item("BSM")
# Row 95
$ws.Range("H94").Value = 2684.2222
$ws.Range("I94").Value = 2684.2222
$ws.Range("K94").Value = 2684.2222
$ws.Range("M94").Value = -2233.2222

# Row 100
$ws.Range("H99").Value = 2270.8572
$ws.Range("I99").Value = 2270.8572
$ws.Range("K99").Value = 2270.8572
$ws.Range("M99").Value = -772.8571999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H16").Value = 453.83334
$ws.Range("I16").Value = 478
$ws.Range("J16").Value = 333
$ws.Range("K16").Value = 478
$ws.Range("L16").Value = 333
$ws.Range("M16").Value = -191
$ws.Range("N16").Value = -907

# Row 100
$ws.Range("H99").Value = 4439.9287
$ws.Range("I99").Value = 4439.9287
$ws.Range("K99").Value = 4439.9287
$ws.Range("M99").Value = -2941.9287

# Row 106
$ws.Range("H105").Value = 26552.5
$ws.Range("I105").Value = 26552.5
$ws.Range("K105").Value = 26552.5
$ws.Range("M105").Value = -24805.5

# Row 114
$ws.Range("H113").Value = 453.83334
$ws.Range("I113").Value = 478
$ws.Range("J113").Value = 333
$ws.Range("K113").Value = 478
$ws.Range("L113").Value = 333
$ws.Range("M113").Value = 1692
$ws.Range("N113").Value = -4673

# Row 123
$ws.Range("H122").Value = 784
$ws.Range("I122").Value = 784
$ws.Range("K122").Value = 2352
$ws.Range("M122").Value = 98

# Row 127
$ws.Range("H126").Value = 4439.9287
$ws.Range("I126").Value = 4439.9287
$ws.Range("K126").Value = 13319.7861
$ws.Range("M126").Value = -10849.7861

# Row 133
$ws.Range("H132").Value = 3596.1
$ws.Range("I132").Value = 3596.1
$ws.Range("K132").Value = 10788.3
$ws.Range("M132").Value = -8258.299999999999

# Row 135
$ws.Range("H134").Value = 1991.1666
$ws.Range("I134").Value = 1886.6
$ws.Range("J134").Value = 2514
$ws.Range("K134").Value = 5659.799999999999
$ws.Range("L134").Value = 7542
$ws.Range("M134").Value = -3124.799999999999
$ws.Range("N134").Value = -12612

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H2").Value = 56.61905
$ws.Range("I2").Value = 21.411764
$ws.Range("J2").Value = 206.25
$ws.Range("K2").Value = 128.470584
$ws.Range("L2").Value = 1237.5
$ws.Range("M2").Value = -15.470584
$ws.Range("N2").Value = -1463.5

# Row 39
$ws.Range("H38").Value = 186.2
$ws.Range("I38").Value = 54.166668
$ws.Range("J38").Value = 384.25
$ws.Range("K38").Value = 162.500004
$ws.Range("L38").Value = 1152.75
$ws.Range("M38").Value = 184.499996
$ws.Range("N38").Value = -1846.75

# Row 61
$ws.Range("H60").Value = 4999.3335
$ws.Range("I60").Value = 4999.3335
$ws.Range("K60").Value = 14998.0005
$ws.Range("M60").Value = -14747.0005

# Row 69
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 72
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 75
$ws.Range("H74").Value = 8000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 78
$ws.Range("H77").Value = 8000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 127
$ws.Range("H126").Value = 4792.8125
$ws.Range("I126").Value = 4075.9092
$ws.Range("K126").Value = 12227.7276
$ws.Range("M126").Value = -9757.7276

# Row 133
$ws.Range("H132").Value = 2498.625
$ws.Range("I132").Value = 2498.5334
$ws.Range("K132").Value = 7495.600199999999
$ws.Range("M132").Value = -4965.600199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 23
$ws.Range("H22").Value = 6639.7
$ws.Range("I22").Value = 2199.5
$ws.Range("J22").Value = 7749.75
$ws.Range("K22").Value = 2199.5
$ws.Range("L22").Value = 7749.75
$ws.Range("M22").Value = -1904.5
$ws.Range("N22").Value = -8339.75

# Row 28
$ws.Range("H27").Value = 6639.7
$ws.Range("I27").Value = 2199.5
$ws.Range("J27").Value = 7749.75
$ws.Range("K27").Value = 2199.5
$ws.Range("L27").Value = 7749.75
$ws.Range("M27").Value = -2092.5
$ws.Range("N27").Value = -7963.75

# Row 47
$ws.Range("H46").Value = 4714.0713
$ws.Range("J46").Value = 5110.778
$ws.Range("L46").Value = 5110.778
$ws.Range("N46").Value = -5486.778

# Row 62
$ws.Range("H61").Value = 809.8
$ws.Range("I61").Value = 787.25
$ws.Range("K61").Value = 787.25
$ws.Range("M61").Value = -585.25

# Row 69
$ws.Range("H68").Value = 2778.7144
$ws.Range("J68").Value = 2483
$ws.Range("L68").Value = 2483
$ws.Range("N68").Value = -3981

# Row 72
$ws.Range("H71").Value = 2778.7144
$ws.Range("J71").Value = 2483
$ws.Range("L71").Value = 12415
$ws.Range("N71").Value = -19903

# Row 114
$ws.Range("H113").Value = 809.8
$ws.Range("I113").Value = 787.25
$ws.Range("K113").Value = 787.25
$ws.Range("M113").Value = 1382.75

# Row 123
$ws.Range("H122").Value = 5542.222
$ws.Range("I122").Value = 4840.2856
$ws.Range("J122").Value = 7999
$ws.Range("K122").Value = 14520.8568
$ws.Range("L122").Value = 23997
$ws.Range("M122").Value = -12070.8568
$ws.Range("N122").Value = -28897

# Row 133
$ws.Range("H132").Value = 1938.7858
$ws.Range("I132").Value = 1933.6923
$ws.Range("K132").Value = 5801.0769
$ws.Range("M132").Value = -3271.0769

# Row 137
$ws.Range("H136").Value = 3491.5
$ws.Range("I136").Value = 3491.5
$ws.Range("K136").Value = 10474.5
$ws.Range("M136").Value = -7924.5

$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H100").Value = 2221.4614
$ws.Range("I100").Value = 1552.9
$ws.Range("K100").Value = 3105.8
$ws.Range("M100").Value = -2564.8
